$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 57: switch heuristics label to "1+2+3.5" (shared string idx 25) and update stats ---
$ws.Range("E57").Value = "1+2+3.5"
$ws.Range("H57").Value = 23
$ws.Range("I57").Value = 174
$ws.Range("J57").Value = 20.15
$ws.Range("K57").Value = 43.83

# --- Row 58: fill in the previously-sparse row with full run parameters ---
$ws.Range("B58").Value = 2
$ws.Range("C58").Value = 10
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = "1+2+3.5"
$ws.Range("F58").Value = 200
$ws.Range("H58").Value = 17
$ws.Range("I58").Value = 182
$ws.Range("J58").Value = 18.62
$ws.Range("K58").Value = 45.32

# --- Row 59 ---
$ws.Range("B59").Value = 3
$ws.Range("C59").Value = 3
$ws.Range("D59").Value = 0
$ws.Range("E59").Value = "1+2+3.5"
$ws.Range("F59").Value = 200
$ws.Range("H59").Value = 19
$ws.Range("I59").Value = 178
$ws.Range("J59").Value = 19.5
$ws.Range("K59").Value = 44.45

# --- Row 60 ---
$ws.Range("B60").Value = 4
$ws.Range("C60").Value = 4
$ws.Range("D60").Value = 0
$ws.Range("E60").Value = "1+2+3.5"
$ws.Range("F60").Value = 100
$ws.Range("H60").Value = 14
$ws.Range("I60").Value = 85
$ws.Range("J60").Value = 20.29
$ws.Range("K60").Value = 43.7

# --- Row 61 ---
$ws.Range("B61").Value = 5
$ws.Range("C61").Value = 5
$ws.Range("D61").Value = 0
$ws.Range("E61").Value = "1+2+3.5"
$ws.Range("F61").Value = 10
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 10
$ws.Range("J61").Value = 17.8
$ws.Range("K61").Value = 46.2

# --- Row 62: brand new row (no column A value) ---
$ws.Range("B62").Value = 5
$ws.Range("C62").Value = 20
$ws.Range("D62").Value = 0
$ws.Range("E62").Value = "1+2+3.5"
$ws.Range("F62").Value = 3
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 3
$ws.Range("J62").Value = 22
$ws.Range("K62").Value = 42

# --- View state: move selection to the newly added row ---
$ws.Range("L62").Select() | Out-Null
